$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C29").Value = 17
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 12
$ws.Range("F29").Value = 0.8605851979345954
